# Daily attendance processing - 2025-10-18 18:47:07
#
# The "Recorded By" column (G) on the "Session Analysis Results" sheet
# lists the user(s) who recorded attendance for each session. A set of
# rows had "dnasr281@gmail.com"/"System"/"admin@admin.com" combinations
# recorded in the wrong order; this corrects the order of the two
# comma-separated entries for exactly those rows, leaving every other
# cell (value, formula, and formatting/style) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The exact "Recorded By" combinations that need their two entries
# swapped back into the correct order.
$fixups = @{
    "dnasr281@gmail.com, System"          = "System, dnasr281@gmail.com"
    "System, admin@admin.com"             = "admin@admin.com, System"
    "dnasr281@gmail.com, admin@admin.com" = "admin@admin.com, dnasr281@gmail.com"
}

$col = 7  # column G = "Recorded By"
$lastRow = 160

for ($row = 1; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, $col)
    $text = $cell.Text

    if ($fixups.ContainsKey($text)) {
        $cell.Value = $fixups[$text]
    }
}
